# "dynamic Prdocut Search and Click"
#
# This edit:
#  1. Renames the shared string "POLO" -> "Beanie" (only referenced by
#     Sheet2!A1).
#  2. Adds a new worksheet "Sheet3" right after "Sheet2" (created via
#     Copy so it inherits Sheet2's formatting/namespace set, then wiped
#     of inherited cell contents).
#  3. Moves the value that used to live in Sheet2!B1 (20) into the new
#     Sheet3!A1 as a plain number, and removes it from Sheet2.
#  4. Updates the selection/active-cell on Sheet2 (no longer the active
#     tab) and makes Sheet3 the active tab with its own selection.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# 1. Product name text change: POLO -> Beanie
$ws2.Range("A1").Value = "Beanie"

# 2. Insert the new "Sheet3" right after "Sheet2" (becomes the last sheet).
#    Copying Sheet2 keeps the worksheet's full namespace/formatting setup.
$ws2.Copy($null, $ws2) | Out-Null
$ws3 = $wb.Worksheets.Item($ws2.Index + 1)
$ws3.Name = "Sheet3"

# 3. Sheet3!A1 gets the plain numeric value 20 (no longer a shared string),
#    and the copied B1 cell is removed so only A1 remains.
$ws3.Range("A1").Value = 20
$ws3.Range("B1").ClearContents()

# Sheet2 no longer holds the moved value
$ws2.Range("B1").ClearContents()

# 4. Update selections: Sheet2 is no longer the active tab, Sheet3 is
$ws2.Activate() | Out-Null
$ws2.Range("D5").Select() | Out-Null

$ws3.Activate() | Out-Null
$ws3.Range("F9").Select() | Out-Null
